# Updated cryptos list refresh: new Price (D) / Volume(1h) (E) figures
# for each coin row, per the source commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    # Write $value into $cell as literal text, even when it looks like a
    # number (e.g. "1.00", "35.60", "0.0₃0911" lookalikes) so Excel does not
    # silently coerce it to a Double and drop the formatting-significant
    # trailing zeros. Number-format is toggled only transiently and the
    # cells original style is restored afterwards so no stray style/
    # number-format ends up attached to the cell.
    param($cell, [string]$value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "41.038.65"
$ws.Range("E2").Value = "  -3.67%  "
Set-TextValue $ws.Range("D3") "2.457.87"
$ws.Range("E3").Value = "  -2.80%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue $ws.Range("D5") "310.54"
$ws.Range("E5").Value = "  +0.29%  "
Set-TextValue $ws.Range("D6") "92.95"
$ws.Range("E6").Value = "  -7.07%  "
$ws.Range("E7").Value = "  -3.00%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -5.34%  "
Set-TextValue $ws.Range("D10") "33.09"
$ws.Range("E10").Value = "  -7.38%  "
$ws.Range("E11").Value = "  -3.61%  "
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("E13").Value = "  -5.70%  "
Set-TextValue $ws.Range("D14") "2.832.18"
$ws.Range("E14").Value = "  -2.97%  "
Set-TextValue $ws.Range("D15") "2.460.16"
$ws.Range("E15").Value = "  -3.06%  "
Set-TextValue $ws.Range("D16") "14.76"
$ws.Range("E16").Value = "  -3.51%  "
Set-TextValue $ws.Range("D17") "0.777"
$ws.Range("E17").Value = "  -3.99%  "
Set-TextValue $ws.Range("D18") "41.001.02"
$ws.Range("E18").Value = "  -3.74%  "
Set-TextValue $ws.Range("D19") "6.25"
$ws.Range("E19").Value = "  -6.73%  "
$ws.Range("E20").Value = "  -3.89%  "
Set-TextValue $ws.Range("D21") "11.03"
$ws.Range("E21").Value = "  -9.67%  "
Set-TextValue $ws.Range("D22") "67.88"
$ws.Range("E22").Value = "  -2.07%  "
Set-TextValue $ws.Range("D23") "234.02"
$ws.Range("E23").Value = "  -3.44%  "
$ws.Range("E24").Value = "  -4.67%  "
$ws.Range("E26").Value = "  -7.60%  "
Set-TextValue $ws.Range("D27") "23.74"
$ws.Range("E27").Value = "  -7.01%  "
$ws.Range("E28").Value = "  -5.85%  "
Set-TextValue $ws.Range("D29") "9.52"
$ws.Range("E29").Value = "  -6.09%  "
Set-TextValue $ws.Range("D30") "35.60"
$ws.Range("E30").Value = "  -7.11%  "
Set-TextValue $ws.Range("D31") "151.01"
$ws.Range("E31").Value = "  -4.09%  "
Set-TextValue $ws.Range("D32") "5.45"
$ws.Range("E32").Value = "  -5.17%  "
Set-TextValue $ws.Range("D33") "2.64"
$ws.Range("E33").Value = "  -5.78%  "
$ws.Range("E34").Value = "  -4.06%  "
$ws.Range("E35").Value = "  -6.15%  "
Set-TextValue $ws.Range("D36") "2.97"
$ws.Range("E36").Value = "  -5.43%  "
Set-TextValue $ws.Range("D37") "1.85"
$ws.Range("E37").Value = "  -5.79%  "
Set-TextValue $ws.Range("D38") "16.60"
$ws.Range("E38").Value = "  -8.53%  "
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("E40").Value = "  -8.45%  "
Set-TextValue $ws.Range("D41") "4.15"
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("E42").Value = "  +0.20%  "
Set-TextValue $ws.Range("D43") "19.85"
$ws.Range("E43").Value = "  -12.08%  "
Set-TextValue $ws.Range("D44") "1.963.49"
Set-TextValue $ws.Range("D45") "0.0282"
$ws.Range("E45").Value = "  -5.90%  "
Set-TextValue $ws.Range("D46") "3.00"
$ws.Range("E46").Value = "  -8.27%  "
Set-TextValue $ws.Range("D47") "8.58"
$ws.Range("E47").Value = "  -3.79%  "
Set-TextValue $ws.Range("D48") "69.18"
$ws.Range("E48").Value = "  -4.11%  "
Set-TextValue $ws.Range("D49") "96.08"
$ws.Range("E49").Value = "  -4.77%  "
$ws.Range("E50").Value = "  -6.84%  "
Set-TextValue $ws.Range("D51") "73.63"
$ws.Range("E51").Value = "  -7.01%  "
